$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 345.7096700033408
$ws.Range("C3").Value = 380.064859628955
$ws.Range("D3").Value = 414.4200492545692
$ws.Range("E3").Value = 448.7752388801832
$ws.Range("F3").Value = 483.1304285057973
$ws.Range("G3").Value = 517.4856181314115
$ws.Range("H3").Value = 551.8408077570257
$ws.Range("I3").Value = 586.1959973826397
$ws.Range("J3").Value = 620.5511870082538
$ws.Range("K3").Value = 654.906376633868
$ws.Range("L3").Value = 689.2615662594819
$ws.Range("M3").Value = 723.6167558850962
$ws.Range("N3").Value = 757.9719455107105
$ws.Range("O3").Value = 792.3271351363244
$ws.Range("P3").Value = 826.6823247619386
$ws.Range("Q3").Value = 861.0375143875527
$ws.Range("R3").Value = 895.3927040131671
$ws.Range("S3").Value = 929.7478936387809
$ws.Range("T3").Value = 964.1030832643951
$ws.Range("U3").Value = 998.4582728900093
$ws.Range("V3").Value = 1032.813462515624
$ws.Range("W3").Value = 1067.168652141238
$ws.Range("X3").Value = 1101.523841766852
$ws.Range("Y3").Value = 1135.879031392466
$ws.Range("Z3").Value = 1170.23422101808
$ws.Range("AA3").Value = 1204.589410643694
$ws.Range("AB3").Value = 1238.944600269308
$ws.Range("AC3").Value = 1273.299789894922
$ws.Range("AD3").Value = 1307.654979520536
$ws.Range("AE3").Value = 1342.01016914615
$ws.Range("AF3").Value = 1376.365358771764
$ws.Range("AG3").Value = 1410.720548397379
$ws.Range("AH3").Value = 1445.075738022993
$ws.Range("AI3").Value = 1479.430927648607
$ws.Range("AJ3").Value = 1513.786117274221
$ws.Range("AK3").Value = 1548.141306899835
$ws.Range("AL3").Value = 1582.496496525449
$ws.Range("AM3").Value = 1616.851686151063
$ws.Range("AN3").Value = 1651.206875776677
$ws.Range("AO3").Value = 1685.562065402292
$ws.Range("AP3").Value = 1719.917255027906
$ws.Range("AQ3").Value = 1754.27244465352
$ws.Range("AR3").Value = 1788.627634279134
$ws.Range("AS3").Value = 1822.982823904748
$ws.Range("AT3").Value = 1857.338013530362
$ws.Range("AU3").Value = 1891.693203155976
$ws.Range("AV3").Value = 1926.04839278159
$ws.Range("AW3").Value = 1960.403582407205
$ws.Range("AX3").Value = 1994.758772032819
$ws.Range("AY3").Value = 2029.113961658433
